# Generate Report for Handback
# The "423069c5-47e9-41c5-a0d5-9576d90d397d.md" file has been successfully
# handed back (it is now in sync with en-US), so every sheet that tracks
# its status/handback information needs to be refreshed to reflect that.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns for the 423069c5 file
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# zh-cn sheet: Status / Latest Handback DateTime / Error Detail for the
# 423069c5 file (row 3)
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Handed back: in sync with en-US"
$zhcn.Range("K3").Value = "2016-08-27 18:55:41"
$zhcn.Range("P3").Value = ""

# ---------------------------------------------------------------------
# de-de sheet: Status / Latest Handback DateTime / Error Detail for the
# 423069c5 file (row 3)
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Handed back: in sync with en-US"
$dede.Range("K3").Value = "2016-08-27 18:55:48"
$dede.Range("P3").Value = ""
